$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 245, shifting existing rows 245:311 down to 246:312
$ws.Rows.Item(245).Insert()

# Populate the newly inserted row 245 with the new weekly record
$ws.Range("A245").Value = 3
$ws.Range("B245").Value = "Femacal de La Calera"
$ws.Range("C245").Value = "Coquimbo"
$ws.Range("D245").Value = 44463
$ws.Range("E245").Value = 5
$ws.Range("F245").Value = 100112045
$ws.Range("G245").Value = "Zapallo"
$ws.Range("H245").Value = "Camote"
$ws.Range("I245").Value = "1a (guarda)"
$ws.Range("J245").Value = 220
$ws.Range("K245").Value = 700
$ws.Range("L245").Value = 750
$ws.Range("M245").Value = 727
$ws.Range("N245").Value = "$/kilo (volumen en unidades)"
$ws.Range("O245").Value = "Provincia de Talca"
$ws.Range("P245").Value = 727
$ws.Range("Q245").Value = 1
$ws.Range("R245").Value = "Hortaliza"
